$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet2")

$ws.Range("E1").Value = "E"
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 10
$ws.Range("E4").Value = 15
